# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C9").Value = 45243
